# Clean up function calls so that aggregate_flows is called properly everywhere,
# and set up the steel sheet for tomorrow's testing.

$wb = $excel.ActiveWorkbook

$wsChains      = $wb.Worksheets.Item("chains")
$wsConnections = $wb.Worksheets.Item("connections")
$wsSteel       = $wb.Worksheets.Item("steel")

# --- "chains" sheet: remove the stray "fuel" chain row (old row 8) ---
$wsChains.Rows("8:8").Delete() | Out-Null

# --- "connections" sheet: remove the obsolete fuel / biofuel connection rows ---
$wsConnections.Rows("19:26").Delete() | Out-Null

# --- "steel" sheet: extend the process chain with the casting / rolling step ---
$wsSteel.Range("A5").Value = "crude steel"
$wsSteel.Range("B5").Value = "simple_casting"
$wsSteel.Range("C5").Value = "hot rolled coil"

# --- Restore/update the view state (selections + active sheet) ---
$wsChains.Activate() | Out-Null
$wsChains.Range("B18").Select() | Out-Null

$wsConnections.Activate() | Out-Null
$wsConnections.Range("A19:XFD26").Select() | Out-Null

$wsSteel.Activate() | Out-Null
$wsSteel.Range("B11").Select() | Out-Null
